$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 31
$ws.Range("D1").Value = 36

$ws.Range("A2").Value = 14
$ws.Range("B2").Value = 26
$ws.Range("C2").Value = 27
$ws.Range("D2").Value = 28
$ws.Range("E2").Value = 34
$ws.Range("F2").Value = 38
$ws.Range("G2").Value = "과거기록 : [564]회차 4등"

$ws.Range("B3").Value = 20
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 39
$ws.Range("F3").Value = 42
$ws.Range("G3").Value = "과거기록 : [349]회차 4등"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 18
$ws.Range("C4").Value = 24
$ws.Range("D4").Value = 26
$ws.Range("E4").Value = 36
$ws.Range("F4").Value = 41
$ws.Range("G4").Value = "과거기록 : [172]회차 4등"

$ws.Range("A5").Value = 10
$ws.Range("B5").Value = 21
$ws.Range("C5").Value = 25
$ws.Range("D5").Value = 29
$ws.Range("E5").Value = 37
$ws.Range("F5").Value = 42
$ws.Range("G5").Value = "과거기록 : [966]회차 4등"
